# Weekly crime data update (new crime data collected)
# Updates the CompStat report header (volume/number + reporting week dates)
# and refreshes the weekly/28-day/YTD/2-year crime statistics table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: bump issue number and advance the reporting week ---
$ws.Range("A8").Value = "Volume 31   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/29/2024  Through  8/4/2024"

# --- Crime statistics table (rows 14-33) ---
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 3
$ws.Range("G14").Value = 9
$ws.Range("H14").Value = 11.111111111111
$ws.Range("I14").Value = 67
$ws.Range("J14").Value = 82
$ws.Range("K14").Value = -18.292682926829
$ws.Range("L14").Value = -23.863636363636
$ws.Range("M14").Value = -15.189873417721
$ws.Range("N14").Value = -77.28813559322
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 75
$ws.Range("G15").Value = 27
$ws.Range("H15").Value = -18.518518518518
$ws.Range("I15").Value = 243
$ws.Range("J15").Value = 231
$ws.Range("K15").Value = 5.194805194805
$ws.Range("L15").Value = -0.816326530612
$ws.Range("M15").Value = 37.28813559322
$ws.Range("N15").Value = -42.688679245283
$ws.Range("D16").Value = 115
$ws.Range("E16").Value = -10.434782608695
$ws.Range("G16").Value = 449
$ws.Range("H16").Value = -4.89977728285
$ws.Range("I16").Value = 2926
$ws.Range("J16").Value = 2812
$ws.Range("K16").Value = 4.054054054054
$ws.Range("L16").Value = -2.205882352941
$ws.Range("M16").Value = 14.70011760094
$ws.Range("N16").Value = -68.91202719932
$ws.Range("C17").Value = 168
$ws.Range("D17").Value = 175
$ws.Range("E17").Value = -4
$ws.Range("F17").Value = 698
$ws.Range("G17").Value = 741
$ws.Range("H17").Value = -5.802968960863
$ws.Range("I17").Value = 4904
$ws.Range("J17").Value = 4793
$ws.Range("K17").Value = 2.315877321093
$ws.Range("L17").Value = 12.091428571428
$ws.Range("M17").Value = 85.82796513831
$ws.Range("N17").Value = -10.133773135422
$ws.Range("C18").Value = 71
$ws.Range("D18").Value = 60
$ws.Range("E18").Value = 18.333333333333
$ws.Range("F18").Value = 250
$ws.Range("G18").Value = 247
$ws.Range("H18").Value = 1.214574898785
$ws.Range("I18").Value = 1733
$ws.Range("J18").Value = 1785
$ws.Range("K18").Value = -2.913165266106
$ws.Range("L18").Value = -0.801373783629
$ws.Range("M18").Value = -7.966011683483
$ws.Range("N18").Value = -84.450426200089
$ws.Range("C19").Value = 193
$ws.Range("D19").Value = 170
$ws.Range("E19").Value = 13.529411764705
$ws.Range("F19").Value = 795
$ws.Range("G19").Value = 634
$ws.Range("H19").Value = 25.394321766561
$ws.Range("I19").Value = 5352
$ws.Range("J19").Value = 4620
$ws.Range("K19").Value = 15.844155844155
$ws.Range("L19").Value = 14.212548015364
$ws.Range("M19").Value = 99.627004848937
$ws.Range("N19").Value = 23.688467760573
$ws.Range("C20").Value = 93
$ws.Range("D20").Value = 107
$ws.Range("E20").Value = -13.084112149532
$ws.Range("F20").Value = 358
$ws.Range("G20").Value = 447
$ws.Range("H20").Value = -19.910514541387
$ws.Range("I20").Value = 2449
$ws.Range("J20").Value = 3172
$ws.Range("K20").Value = -22.793190416141
$ws.Range("L20").Value = 4.524114383269
$ws.Range("M20").Value = 100.90237899918
$ws.Range("N20").Value = -72.930253122582
$ws.Range("C21").Value = 638
$ws.Range("D21").Value = 634
$ws.Range("E21").Value = 0.630914826498
$ws.Range("F21").Value = 2560
$ws.Range("G21").Value = 2554
$ws.Range("H21").Value = 0.234925606891
$ws.Range("I21").Value = 17674
$ws.Range("J21").Value = 17495
$ws.Range("K21").Value = 1.023149471277
$ws.Range("L21").Value = 7.271182325807
$ws.Range("M21").Value = 57.396028141419
$ws.Range("N21").Value = -55.932879547211
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 166.666666666667
$ws.Range("G22").Value = 19
$ws.Range("H22").Value = 15.78947368421
$ws.Range("I22").Value = 195
$ws.Range("J22").Value = 173
$ws.Range("K22").Value = 12.71676300578
$ws.Range("L22").Value = -7.142857142857
$ws.Range("M22").Value = 0.515463917525
$ws.Range("C23").Value = 41
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = 28.125
$ws.Range("F23").Value = 139
$ws.Range("G23").Value = 148
$ws.Range("H23").Value = -6.081081081081
$ws.Range("I23").Value = 1008
$ws.Range("J23").Value = 1058
$ws.Range("K23").Value = -4.725897920604
$ws.Range("L23").Value = 5.882352941176
$ws.Range("M23").Value = 59.241706161137
$ws.Range("C24").Value = 369
$ws.Range("D24").Value = 344
$ws.Range("E24").Value = 7.267441860465
$ws.Range("F24").Value = 1393
$ws.Range("G24").Value = 1499
$ws.Range("H24").Value = -7.071380920613
$ws.Range("I24").Value = 9569
$ws.Range("J24").Value = 10620
$ws.Range("K24").Value = -9.896421845574
$ws.Range("L24").Value = -12.283435695297
$ws.Range("M24").Value = 28.391251844894
$ws.Range("C25").Value = 164
$ws.Range("D25").Value = 125
$ws.Range("E25").Value = 31.2
$ws.Range("F25").Value = 572
$ws.Range("G25").Value = 606
$ws.Range("H25").Value = -5.610561056105
$ws.Range("I25").Value = 3836
$ws.Range("J25").Value = 4553
$ws.Range("K25").Value = -15.747858554799
$ws.Range("L25").Value = -31.949618591449
$ws.Range("C26").Value = 225
$ws.Range("D26").Value = 209
$ws.Range("E26").Value = 7.655502392344
$ws.Range("F26").Value = 926
$ws.Range("G26").Value = 834
$ws.Range("H26").Value = 11.031175059952
$ws.Range("I26").Value = 6649
$ws.Range("J26").Value = 6244
$ws.Range("K26").Value = 6.486226777706
$ws.Range("L26").Value = 10.430161102806
$ws.Range("M26").Value = 0.211002260738
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = 14
$ws.Range("E27").Value = -28.571428571428
$ws.Range("F27").Value = 39
$ws.Range("G27").Value = 50
$ws.Range("H27").Value = -22
$ws.Range("I27").Value = 383
$ws.Range("J27").Value = 389
$ws.Range("K27").Value = -1.542416452442
$ws.Range("L27").Value = -8.809523809523
$ws.Range("C28").Value = 21
$ws.Range("D28").Value = 11
$ws.Range("E28").Value = 90.90909090909
$ws.Range("F28").Value = 97
$ws.Range("G28").Value = 74
$ws.Range("H28").Value = 31.081081081081
$ws.Range("I28").Value = 726
$ws.Range("J28").Value = 618
$ws.Range("K28").Value = 17.475728155339
$ws.Range("L28").Value = 36.466165413533
$ws.Range("C29").Value = 9
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 80
$ws.Range("F29").Value = 42
$ws.Range("G29").Value = 32
$ws.Range("H29").Value = 31.25
$ws.Range("I29").Value = 242
$ws.Range("J29").Value = 239
$ws.Range("K29").Value = 1.255230125523
$ws.Range("L29").Value = -23.417721518987
$ws.Range("M29").Value = -11.678832116788
$ws.Range("N29").Value = -70.666666666666
$ws.Range("C30").Value = 7
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = 40
$ws.Range("F30").Value = 32
$ws.Range("G30").Value = 24
$ws.Range("H30").Value = 33.333333333333
$ws.Range("I30").Value = 195
$ws.Range("J30").Value = 194
$ws.Range("K30").Value = 0.515463917525
$ws.Range("L30").Value = -27.777777777777
$ws.Range("M30").Value = -15.217391304347
$ws.Range("N30").Value = -74.10358565737
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 13
$ws.Range("K31").Value = 7.692307692307
$ws.Range("C33").Value = 1
$ws.Range("I33").Value = 27
$ws.Range("K33").Value = 17.391304347826
$ws.Range("L33").Value = -37.209302325581

# C33 (Traffic Fatalities / Week to Date 2024) was blank ("0" placeholder
# text) and is now a real count; match the numeric formatting already used
# by the other count cells in this table (e.g. F33).
$ws.Range("C33").NumberFormat = $ws.Range("F33").NumberFormat
